$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the style of an existing header cell (e.g. AC1) to the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats

# Data rows 2-70: Wins=62, Losses=100, Ties=0
for ($r = 2; $r -le 70; $r++) {
    $ws.Cells.Item($r, 30).Value = 62   # AD
    $ws.Cells.Item($r, 31).Value = 100  # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
